$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44301
$ws.Range("M2").Value = 38

# Row 3
$ws.Range("D3").Value = 44377
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44298
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range("S4").Value = 1100

# Row 5
$ws.Range("D5").Value = 44300
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 22000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 22000
$ws.Range("S5").Value = 1100

# Row 6
$ws.Range("D6").Value = 44403
$ws.Range("M6").Value = 50

# Row 8
$ws.Range("D8").Value = 44385
$ws.Range("M8").Value = 36
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44406
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 1000

# Row 11
$ws.Range("D11").Value = 44305
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 22000
$ws.Range("P11").Value = 22000
$ws.Range("S11").Value = 1100

# Row 12
$ws.Range("D12").Value = 44445
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("S12").Value = 1000

# Row 13
$ws.Range("D13").Value = 44307
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 22000
$ws.Range("O13").Value = 22000
$ws.Range("P13").Value = 22000
$ws.Range("S13").Value = 1100

# Row 15
$ws.Range("D15").Value = 44376
$ws.Range("M15").Value = 38
$ws.Range("N15").Value = 20000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 20000
$ws.Range("S15").Value = 1000

# Row 16
$ws.Range("D16").Value = 44294
$ws.Range("M16").Value = 25
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 25000
$ws.Range("P16").Value = 25000
$ws.Range("S16").Value = 1250

# Row 17
$ws.Range("D17").Value = 44389

# Row 19
$ws.Range("D19").Value = 44292
$ws.Range("M19").Value = 30
$ws.Range("N19").Value = 25000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 25000
$ws.Range("S19").Value = 1250
